$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of the
# existing header cells (H1 uses style index 1: bold/border/center).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I and J, rows 2-33: row number -> (I value, J value)
$values = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 6)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(1, 6)
    17 = @(1, 5)
    18 = @(1, 5)
    19 = @(1, 4)
    20 = @(9, 9)
    21 = @(3, 4)
    22 = @(6, 7)
    23 = @(6, 7)
    24 = @(8, 8)
    25 = @(8, 9)
    26 = @(7, 7)
    27 = @(3, 5)
    28 = @(7, 8)
    29 = @(5, 6)
    30 = @(8, 9)
    31 = @(5, 5)
    32 = @(6, 7)
    33 = @(3, 4)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Output "Added I0/IF columns for rows 1-33"
